$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.169.59"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.59%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.736.89"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.74%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "617.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "184.83"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.42%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.736.44"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.44%  "
$ws.Range("E8").Value = "  -3.79%  "
$ws.Range("E9").Value = "  -0.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.723"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.00%  "
$ws.Range("E11").Value = "  -7.43%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "57.39"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.91%  "
$ws.Range("E13").Value = "  -6.03%  "
$ws.Range("E14").Value = "  -5.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.332.15"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.02%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.743.06"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.52"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.08"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.53%  "
$ws.Range("E19").Value = "  -1.58%  "
$ws.Range("E20").Value = "  -4.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "69.040.13"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "415.16"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.57%  "
$ws.Range("E23").Value = "  -1.61%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "89.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.62%  "
$ws.Range("E25").Value = "  -4.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.85"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.29%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.02"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.97%  "
$ws.Range("E29").Value = "  -3.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.68"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.69%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.34"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.01%  "
$ws.Range("E32").Value = "  -15.44%  "
$ws.Range("E33").Value = "  -4.98%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.123"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.38%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "44.75"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "621.14"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "65.90"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.07%  "
$ws.Range("E38").Value = "  -8.68%  "
$ws.Range("B39").Value = "Dai"
$ws.Range("C39").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.17%  "
$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.406"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.73%  "
$ws.Range("E41").Value = "  -0.15%  "
$ws.Range("E42").Value = "  -0.24%  "
$ws.Range("E43").Value = "  -4.94%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0446"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.65"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.28%  "
$ws.Range("B46").Value = "THORChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.29"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.33%  "
$ws.Range("E47").Value = "  -2.67%  "
$ws.Range("B48").Value = "dogwifhat"
$ws.Range("C48").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.78"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -17.64%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.827.66"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.58%  "
$ws.Range("E50").Value = "  -3.55%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.11"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.35%  "
